$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Append three new rows of data to the table (rows 28-30), mirroring the
# existing rows' layout: B=Plaque, C=Type, D=Marque, E=Modele, F=Annee,
# G=Nom, H=Choix, I=Descriptif.

# Copy formatting from the previous last row (27) down to the new rows so
# the thick-bottom border / fills / alignment match the existing table.
$ws.Range("B27:I27").Copy() | Out-Null
$ws.Range("B28:I30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 28 - a first new car entry.
$ws.Range("B28").Value = "300 SL 1957"
$ws.Range("C28").Value = "Voiture"
$ws.Range("D28").Value = "Mercedes"
$ws.Range("E28").Value = "Benz"
$ws.Range("F28").Value = 1970
$ws.Range("G28").Value = "Mathias"
$ws.Range("H28").Value = "OK"
$ws.Range("I28").Value = "RAS"

# Jot down the plate numbers for the next two entries first.
$ws.Range("B29").Value = "BP-931-LB"
$ws.Range("B30").Value = "AM-666-EE"

# Then fill in the rest of row 29.
$ws.Range("C29").Value = "Voiture"
$ws.Range("D29").Value = "BMW"
$ws.Range("E29").Value = "Truc"
$ws.Range("F29").Value = 2005
$ws.Range("G29").Value = "Nicolas"
$ws.Range("H29").Value = "OK"
$ws.Range("I29").Value = "RAS"

# And finally the rest of row 30.
$ws.Range("C30").Value = "Voiture"
$ws.Range("D30").Value = "Citroen"
$ws.Range("E30").Value = "Truc"
$ws.Range("F30").Value = 2010
$ws.Range("G30").Value = "Kaelig"
$ws.Range("H30").Value = "OK"
$ws.Range("I30").Value = "RAS"

$ws.Range("B30").Select() | Out-Null
